# Update the "Variables List Indices" lower-right cell reference from E177 to E178
# (new delivery variables pushed the data range down one row; aggregated variables
# must be converted to TAF).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = "E178"

# Reflect the updated view/selection state captured in the workbook: the window is
# scrolled so column B is the left-most visible column, and D13 (the cell we just
# edited) is the active selection.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D13").Select()
